$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells receiving new values that look like plain numbers (single decimal point)
# must be pre-formatted as Text so Excel keeps them as strings (matching the
# source workbook, which stores these as text, e.g. "306.11" not 306.11).
$textCells = @("D5", "D6", "D10", "D11", "D15", "D17", "D19", "D20", "D22", "D23", "D25", "D27", "D30", "D31", "D32", "D35", "D37", "D40", "D43", "D45", "D46", "D47")
foreach ($cell in $textCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range("D2").Value = "42.947.73"
$ws.Range("E2").Value = "  -0.81%  "
$ws.Range("D3").Value = "2.332.37"
$ws.Range("E3").Value = "  +1.25%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "306.11"
$ws.Range("E5").Value = "  -1.39%  "
$ws.Range("D6").Value = "100.35"
$ws.Range("E6").Value = "  -2.57%  "
$ws.Range("E7").Value = "  -3.78%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("E9").Value = "  -4.10%  "
$ws.Range("D10").Value = "34.73"
$ws.Range("E10").Value = "  -2.56%  "
$ws.Range("D11").Value = "52.33"
$ws.Range("E11").Value = "  +0.75%  "
$ws.Range("E12").Value = "  -1.79%  "
$ws.Range("E13").Value = "  +0.79%  "
$ws.Range("E14").Value = "  -2.23%  "
$ws.Range("D15").Value = "15.87"
$ws.Range("E15").Value = "  +6.13%  "
$ws.Range("D16").Value = "2.310.85"
$ws.Range("E16").Value = "  -2.13%  "
$ws.Range("D17").Value = "0.812"
$ws.Range("E17").Value = "  +0.81%  "
$ws.Range("D18").Value = "42.882.99"
$ws.Range("D19").Value = "11.75"
$ws.Range("E19").Value = "  -4.21%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "6.20"
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "0.0₃0910"
$ws.Range("E21").Value = "  -2.31%  "
$ws.Range("D22").Value = "67.84"
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("D23").Value = "236.35"
$ws.Range("E23").Value = "  -2.12%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").Value = "2.56"
$ws.Range("E25").Value = "  -2.11%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").Value = "25.33"
$ws.Range("E27").Value = "  +1.66%  "
$ws.Range("E28").Value = "  -0.74%  "
$ws.Range("E29").Value = "  +0.96%  "
$ws.Range("D30").Value = "35.04"
$ws.Range("E30").Value = "  -4.07%  "
$ws.Range("D31").Value = "9.39"
$ws.Range("E31").Value = "  -2.55%  "
$ws.Range("D32").Value = "163.48"
$ws.Range("E32").Value = "  -3.95%  "
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("E34").Value = "  -3.21%  "
$ws.Range("D35").Value = "17.61"
$ws.Range("E35").Value = "  -0.77%  "
$ws.Range("E36").Value = "  -3.89%  "
$ws.Range("D37").Value = "4.58"
$ws.Range("E37").Value = "  +5.26%  "
$ws.Range("E38").Value = "  -2.12%  "
$ws.Range("E39").Value = "  -1.44%  "
$ws.Range("D40").Value = "2.91"
$ws.Range("E40").Value = "  -4.85%  "
$ws.Range("E41").Value = "  -3.96%  "
$ws.Range("E42").Value = "  -2.33%  "
$ws.Range("D43").Value = "2.63"
$ws.Range("E43").Value = "  +14.50%  "
$ws.Range("D44").Value = "2.006.42"
$ws.Range("E44").Value = "  +2.09%  "
$ws.Range("D45").Value = "0.0285"
$ws.Range("E45").Value = "  -1.55%  "
$ws.Range("D46").Value = "18.69"
$ws.Range("E46").Value = "  -2.38%  "
$ws.Range("D47").Value = "10.14"
$ws.Range("E47").Value = "  +2.07%  "
$ws.Range("E48").Value = "  -2.62%  "
$ws.Range("E49").Value = "  +0.62%  "
$ws.Range("E50").Value = "  -0.57%  "
$ws.Range("D51").Value = "2.558.94"
$ws.Range("E51").Value = "  +1.16%  "

# Restore default (General) number format now that the text values are set,
# so formatting stays as close as possible to the original workbook.
foreach ($cell in $textCells) {
    $ws.Range($cell).NumberFormat = "General"
}
